$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1:T21").ClearContents()

# Row 1
$ws.Range("A1").Value = "אגמים"
$ws.Range("B1").Value = "תולעת"
$ws.Range("C1").Value = "רעמים"
$ws.Range("D1").Value = "קליפה"
$ws.Range("E1").Value = "חרקים"
$ws.Range("F1").Value = "פרעוש"
$ws.Range("G1").Value = "אריות"
$ws.Range("H1").Value = "תותים"
$ws.Range("I1").Value = "בהמות"
$ws.Range("J1").Value = "גשמים"
$ws.Range("K1").Value = "עננים"
$ws.Range("L1").Value = "בצלים"
$ws.Range("M1").Value = "סופות"
$ws.Range("N1").Value = "חולות"
$ws.Range("O1").Value = "גייזר"
$ws.Range("P1").Value = "חולדה"
$ws.Range("Q1").Value = "זחלים"
$ws.Range("R1").Value = "קוקוס"
$ws.Range("S1").Value = "בננות"
$ws.Range("T1").Value = "פילים"

# Row 2
$ws.Range("A2").Value = "מסגרת"
$ws.Range("C2").Value = "מסגרת"
$ws.Range("D2").Value = "מסגרת"
$ws.Range("E2").Value = "מסגרת"
$ws.Range("F2").Value = "מסגרת"
$ws.Range("H2").Value = "מסגרת"
$ws.Range("J2").Value = "מסגרת"
$ws.Range("K2").Value = "מסגרת"
$ws.Range("L2").Value = "מסגרת"
$ws.Range("O2").Value = "מסגרת"
$ws.Range("P2").Value = "מסגרת"
$ws.Range("Q2").Value = "מסגרת"
$ws.Range("R2").Value = "מסגרת"
$ws.Range("T2").Value = "מסגרת"

# Row 3
$ws.Range("A3").Value = "מלתחה"
$ws.Range("B3").Value = "מלתחה"
$ws.Range("C3").Value = "מלתחה"
$ws.Range("E3").Value = "מלתחה"
$ws.Range("F3").Value = "מלתחה"
$ws.Range("G3").Value = "מלתחה"
$ws.Range("I3").Value = "מלתחה"
$ws.Range("J3").Value = "מלתחה"
$ws.Range("K3").Value = "מלתחה"
$ws.Range("L3").Value = "מלתחה"
$ws.Range("M3").Value = "מלתחה"
$ws.Range("N3").Value = "מלתחה"
$ws.Range("O3").Value = "מלתחה"
$ws.Range("Q3").Value = "מלתחה"
$ws.Range("R3").Value = "מלתחה"
$ws.Range("S3").Value = "מלתחה"
$ws.Range("T3").Value = "מלתחה"

# Row 4
$ws.Range("A4").Value = "שטרות"
$ws.Range("C4").Value = "שטרות"
$ws.Range("D4").Value = "שטרות"
$ws.Range("E4").Value = "שטרות"
$ws.Range("H4").Value = "שטרות"
$ws.Range("J4").Value = "שטרות"
$ws.Range("K4").Value = "שטרות"
$ws.Range("L4").Value = "שטרות"
$ws.Range("O4").Value = "שטרות"
$ws.Range("P4").Value = "שטרות"
$ws.Range("Q4").Value = "שטרות"
$ws.Range("T4").Value = "שטרות"

# Row 5
$ws.Range("B5").Value = "אגורה"
$ws.Range("C5").Value = "אגורה"
$ws.Range("E5").Value = "אגורה"
$ws.Range("F5").Value = "אגורה"
$ws.Range("H5").Value = "אגורה"
$ws.Range("I5").Value = "אגורה"
$ws.Range("J5").Value = "אגורה"
$ws.Range("K5").Value = "אגורה"
$ws.Range("L5").Value = "אגורה"
$ws.Range("M5").Value = "אגורה"
$ws.Range("N5").Value = "אגורה"
$ws.Range("O5").Value = "אגורה"
$ws.Range("Q5").Value = "אגורה"
$ws.Range("R5").Value = "אגורה"
$ws.Range("S5").Value = "אגורה"
$ws.Range("T5").Value = "אגורה"

# Row 6
$ws.Range("A6").Value = "מבחנה"
$ws.Range("B6").Value = "מבחנה"
$ws.Range("C6").Value = "מבחנה"
$ws.Range("E6").Value = "מבחנה"
$ws.Range("F6").Value = "מבחנה"
$ws.Range("G6").Value = "מבחנה"
$ws.Range("H6").Value = "מבחנה"
$ws.Range("I6").Value = "מבחנה"
$ws.Range("J6").Value = "מבחנה"
$ws.Range("K6").Value = "מבחנה"
$ws.Range("L6").Value = "מבחנה"
$ws.Range("M6").Value = "מבחנה"
$ws.Range("N6").Value = "מבחנה"
$ws.Range("O6").Value = "מבחנה"
$ws.Range("Q6").Value = "מבחנה"
$ws.Range("R6").Value = "מבחנה"
$ws.Range("S6").Value = "מבחנה"
$ws.Range("T6").Value = "מבחנה"

# Row 7
$ws.Range("A7").Value = "צינור"
$ws.Range("B7").Value = "צינור"
$ws.Range("C7").Value = "צינור"
$ws.Range("D7").Value = "צינור"
$ws.Range("E7").Value = "צינור"
$ws.Range("H7").Value = "צינור"
$ws.Range("J7").Value = "צינור"
$ws.Range("L7").Value = "צינור"
$ws.Range("P7").Value = "צינור"
$ws.Range("Q7").Value = "צינור"

# Row 8
$ws.Range("A8").Value = "רובוט"
$ws.Range("D8").Value = "רובוט"
$ws.Range("E8").Value = "רובוט"
$ws.Range("J8").Value = "רובוט"
$ws.Range("K8").Value = "רובוט"
$ws.Range("L8").Value = "רובוט"
$ws.Range("O8").Value = "רובוט"
$ws.Range("Q8").Value = "רובוט"
$ws.Range("T8").Value = "רובוט"

# Row 9
$ws.Range("A9").Value = "מכונה"
$ws.Range("B9").Value = "מכונה"
$ws.Range("C9").Value = "מכונה"
$ws.Range("E9").Value = "מכונה"
$ws.Range("F9").Value = "מכונה"
$ws.Range("G9").Value = "מכונה"
$ws.Range("H9").Value = "מכונה"
$ws.Range("I9").Value = "מכונה"
$ws.Range("J9").Value = "מכונה"
$ws.Range("K9").Value = "מכונה"
$ws.Range("L9").Value = "מכונה"
$ws.Range("M9").Value = "מכונה"
$ws.Range("N9").Value = "מכונה"
$ws.Range("O9").Value = "מכונה"
$ws.Range("Q9").Value = "מכונה"
$ws.Range("R9").Value = "מכונה"
$ws.Range("S9").Value = "מכונה"
$ws.Range("T9").Value = "מכונה"

# Row 10
$ws.Range("B10").Value = "בגדים"
$ws.Range("D10").Value = "בגדים"
$ws.Range("F10").Value = "בגדים"
$ws.Range("G10").Value = "בגדים"
$ws.Range("M10").Value = "בגדים"
$ws.Range("N10").Value = "בגדים"
$ws.Range("O10").Value = "בגדים"
$ws.Range("P10").Value = "בגדים"
$ws.Range("R10").Value = "בגדים"

# Row 11
$ws.Range("A11").Value = "קופסה"
$ws.Range("C11").Value = "קופסה"
$ws.Range("E11").Value = "קופסה"
$ws.Range("F11").Value = "קופסה"
$ws.Range("G11").Value = "קופסה"
$ws.Range("I11").Value = "קופסה"
$ws.Range("J11").Value = "קופסה"
$ws.Range("K11").Value = "קופסה"
$ws.Range("L11").Value = "קופסה"
$ws.Range("O11").Value = "קופסה"
$ws.Range("Q11").Value = "קופסה"
$ws.Range("S11").Value = "קופסה"
$ws.Range("T11").Value = "קופסה"

# Row 12
$ws.Range("A12").Value = "משקפת"
$ws.Range("C12").Value = "משקפת"
$ws.Range("F12").Value = "משקפת"
$ws.Range("H12").Value = "משקפת"
$ws.Range("K12").Value = "משקפת"
$ws.Range("L12").Value = "משקפת"
$ws.Range("O12").Value = "משקפת"
$ws.Range("P12").Value = "משקפת"
$ws.Range("Q12").Value = "משקפת"
$ws.Range("T12").Value = "משקפת"

# Row 13
$ws.Range("A13").Value = "וודקה"
$ws.Range("C13").Value = "וודקה"
$ws.Range("E13").Value = "וודקה"
$ws.Range("F13").Value = "וודקה"
$ws.Range("G13").Value = "וודקה"
$ws.Range("I13").Value = "וודקה"
$ws.Range("J13").Value = "וודקה"
$ws.Range("K13").Value = "וודקה"
$ws.Range("L13").Value = "וודקה"
$ws.Range("O13").Value = "וודקה"
$ws.Range("Q13").Value = "וודקה"
$ws.Range("S13").Value = "וודקה"
$ws.Range("T13").Value = "וודקה"

# Row 14
$ws.Range("B14").Value = "חגורה"
$ws.Range("C14").Value = "חגורה"
$ws.Range("F14").Value = "חגורה"
$ws.Range("G14").Value = "חגורה"
$ws.Range("H14").Value = "חגורה"
$ws.Range("I14").Value = "חגורה"
$ws.Range("J14").Value = "חגורה"
$ws.Range("K14").Value = "חגורה"
$ws.Range("L14").Value = "חגורה"
$ws.Range("M14").Value = "חגורה"
$ws.Range("O14").Value = "חגורה"
$ws.Range("Q14").Value = "חגורה"
$ws.Range("R14").Value = "חגורה"
$ws.Range("S14").Value = "חגורה"
$ws.Range("T14").Value = "חגורה"

# Row 15
$ws.Range("A15").Value = "חולצה"
$ws.Range("C15").Value = "חולצה"
$ws.Range("F15").Value = "חולצה"
$ws.Range("G15").Value = "חולצה"
$ws.Range("I15").Value = "חולצה"
$ws.Range("J15").Value = "חולצה"
$ws.Range("K15").Value = "חולצה"
$ws.Range("O15").Value = "חולצה"
$ws.Range("S15").Value = "חולצה"

# Row 16
$ws.Range("A16").Value = "שבשבת"
$ws.Range("C16").Value = "שבשבת"
$ws.Range("D16").Value = "שבשבת"
$ws.Range("E16").Value = "שבשבת"
$ws.Range("F16").Value = "שבשבת"
$ws.Range("H16").Value = "שבשבת"
$ws.Range("J16").Value = "שבשבת"
$ws.Range("K16").Value = "שבשבת"
$ws.Range("L16").Value = "שבשבת"
$ws.Range("O16").Value = "שבשבת"
$ws.Range("P16").Value = "שבשבת"
$ws.Range("Q16").Value = "שבשבת"
$ws.Range("R16").Value = "שבשבת"
$ws.Range("T16").Value = "שבשבת"

# Row 17
$ws.Range("A17").Value = "חיתול"
$ws.Range("B17").Value = "חיתול"
$ws.Range("C17").Value = "חיתול"
$ws.Range("D17").Value = "חיתול"
$ws.Range("J17").Value = "חיתול"
$ws.Range("K17").Value = "חיתול"
$ws.Range("L17").Value = "חיתול"
$ws.Range("Q17").Value = "חיתול"

# Row 18
$ws.Range("A18").Value = "מקלחת"
$ws.Range("C18").Value = "מקלחת"
$ws.Range("D18").Value = "מקלחת"
$ws.Range("E18").Value = "מקלחת"
$ws.Range("F18").Value = "מקלחת"
$ws.Range("H18").Value = "מקלחת"
$ws.Range("J18").Value = "מקלחת"
$ws.Range("K18").Value = "מקלחת"
$ws.Range("O18").Value = "מקלחת"
$ws.Range("R18").Value = "מקלחת"

# Row 19
$ws.Range("A19").Value = "שמיכה"
$ws.Range("B19").Value = "שמיכה"
$ws.Range("C19").Value = "שמיכה"
$ws.Range("E19").Value = "שמיכה"
$ws.Range("F19").Value = "שמיכה"
$ws.Range("H19").Value = "שמיכה"
$ws.Range("I19").Value = "שמיכה"
$ws.Range("J19").Value = "שמיכה"
$ws.Range("K19").Value = "שמיכה"
$ws.Range("L19").Value = "שמיכה"
$ws.Range("M19").Value = "שמיכה"
$ws.Range("N19").Value = "שמיכה"
$ws.Range("Q19").Value = "שמיכה"
$ws.Range("R19").Value = "שמיכה"
$ws.Range("S19").Value = "שמיכה"
$ws.Range("T19").Value = "שמיכה"

# Row 20
$ws.Range("B20").Value = "שרביט"
$ws.Range("D20").Value = "שרביט"
$ws.Range("I20").Value = "שרביט"
$ws.Range("M20").Value = "שרביט"
$ws.Range("N20").Value = "שרביט"
$ws.Range("O20").Value = "שרביט"
$ws.Range("P20").Value = "שרביט"
$ws.Range("R20").Value = "שרביט"
$ws.Range("S20").Value = "שרביט"

# Row 21
$ws.Range("B21").Value = "מסטיק"
$ws.Range("D21").Value = "מסטיק"
$ws.Range("F21").Value = "מסטיק"
$ws.Range("G21").Value = "מסטיק"
$ws.Range("I21").Value = "מסטיק"
$ws.Range("M21").Value = "מסטיק"
$ws.Range("N21").Value = "מסטיק"
$ws.Range("O21").Value = "מסטיק"
$ws.Range("P21").Value = "מסטיק"
$ws.Range("R21").Value = "מסטיק"
$ws.Range("S21").Value = "מסטיק"
